# Refresh the crypto price/volume table: column D ("Price") and column E
# ("Volume(1h)") are updated with the latest scraped figures.
#
# Column D cells hold plain text (e.g. 27.514.30, 1.004), not numbers. Setting
# Range.Value to a numeric-looking string makes Excel COM automatically store it
# as a real number instead (silently dropping things such as trailing zeros), which
# would corrupt the intended text content. A leading apostrophe forces Excel to keep
# the assignment as text, the same as manually typing an apostrophe before a number
# in a cell; the apostrophe itself is not stored as part of the cell value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.503.09"
$ws.Range("E2").Value = "  +5.44%  "
$ws.Range("D3").Value = "'1.725.89"
$ws.Range("E3").Value = "  +4.74%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'225.58"
$ws.Range("E5").Value = "  +3.43%  "
$ws.Range("D6").Value = "'0.5346"
$ws.Range("E6").Value = "  +2.86%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.2662"
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").Value = "'0.06585"
$ws.Range("E9").Value = "  +4.53%  "
$ws.Range("D10").Value = "'21.59"
$ws.Range("E10").Value = "  +6.32%  "
$ws.Range("D11").Value = "'0.07706"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "'4.604"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "'1.726.12"
$ws.Range("E13").Value = "  +5.77%  "
$ws.Range("D14").Value = "'1.963.82"
$ws.Range("E14").Value = "  +4.73%  "
$ws.Range("D15").Value = "'0.5814"
$ws.Range("E15").Value = "  +4.40%  "
$ws.Range("D16").Value = "'0.0₅8274"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("D17").Value = "'67.78"
$ws.Range("E17").Value = "  +4.13%  "
$ws.Range("D18").Value = "'27.510.24"
$ws.Range("E18").Value = "  +5.54%  "
$ws.Range("D19").Value = "'217.94"
$ws.Range("E19").Value = "  +13.07%  "
$ws.Range("D21").Value = "'4.724"
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("E22").Value = "  +1.70%  "
$ws.Range("D23").Value = "'6.078"
$ws.Range("E23").Value = "  +2.88%  "
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").Value = "'143.84"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "'1.755"
$ws.Range("E26").Value = "  +16.11%  "
$ws.Range("D27").Value = "'0.1236"
$ws.Range("E27").Value = "  +4.66%  "
$ws.Range("D28").Value = "'7.394"
$ws.Range("E28").Value = "  +3.02%  "
$ws.Range("D29").Value = "'16.52"
$ws.Range("E29").Value = "  +4.42%  "
$ws.Range("D30").Value = "'0.05491"
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("D31").Value = "'1.302"
$ws.Range("E31").Value = "  +2.64%  "
$ws.Range("D32").Value = "'3.561"
$ws.Range("E32").Value = "  +3.43%  "
$ws.Range("D33").Value = "'3.438"
$ws.Range("E33").Value = "  +3.48%  "
$ws.Range("D34").Value = "'1.655"
$ws.Range("E34").Value = "  +6.91%  "
$ws.Range("D35").Value = "'2.861"
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("D36").Value = "'0.9639"
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("D37").Value = "'2.426"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").Value = "'0.5970"
$ws.Range("E38").Value = "  +6.82%  "
$ws.Range("D39").Value = "'0.01647"
$ws.Range("E39").Value = "  +4.82%  "
$ws.Range("D40").Value = "'5.903"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("D41").Value = "'1.054.24"
$ws.Range("E41").Value = "  +2.72%  "
$ws.Range("D42").Value = "'0.8523"
$ws.Range("D44").Value = "'101.39"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").Value = "'1.869.97"
$ws.Range("D46").Value = "'0.0₈117"
$ws.Range("E46").Value = "  +3.82%  "
$ws.Range("D47").Value = "'58.87"
$ws.Range("E47").Value = "  +2.78%  "
$ws.Range("D48").Value = "'0.4475"
$ws.Range("E48").Value = "  +3.70%  "
$ws.Range("D49").Value = "'8.216"
$ws.Range("E49").Value = "  +3.95%  "
$ws.Range("D50").Value = "'1.002"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "'0.05245"
$ws.Range("E51").Value = "  +2.63%  "
